$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Shift the two existing week blocks back by 2 days
#    (B2:B8 was 43599..43605 -> 43597..43603)
#    (B11:B17 was 43606..43612 -> 43604..43610)
# ------------------------------------------------------------------
for ($r = 2; $r -le 8; $r++) {
  $cell = $ws.Cells.Item($r, 2)
  $cell.Value = $cell.Value2 - 2
}
for ($r = 11; $r -le 17; $r++) {
  $cell = $ws.Cells.Item($r, 2)
  $cell.Value = $cell.Value2 - 2
}

# ------------------------------------------------------------------
# 2) Build two more week blocks at rows 19-27 and 28-36, matching the
#    layout of the existing ones: a header row (A:G, text), seven day
#    rows (A = day name, B = date only), and a styled blank trailer
#    row (A:G, style of row 9/18).
# ------------------------------------------------------------------
$ws.Range("A1:G1").Copy($ws.Range("A19"))
$ws.Range("A9:G9").Copy($ws.Range("A27"))

$ws.Range("A1:G1").Copy($ws.Range("A28"))
$ws.Range("A9:G9").Copy($ws.Range("A36"))

$days = @("Sun", "Mon", "Tue", "Wed", "Thu", "Fri", "Sat")

# block 3: rows 20-26, dates 43611..43617
$startDate3 = 43611
for ($i = 0; $i -le 6; $i++) {
  $r = 20 + $i
  $ws.Cells.Item($r, 1).Value = $days[$i]
  $ws.Cells.Item($r, 2).Value = $startDate3 + $i
}

# block 4: rows 29-35, dates 43618..43624
$startDate4 = 43618
for ($i = 0; $i -le 6; $i++) {
  $r = 29 + $i
  $ws.Cells.Item($r, 1).Value = $days[$i]
  $ws.Cells.Item($r, 2).Value = $startDate4 + $i
}

# ------------------------------------------------------------------
# 3) Friday of the new 3rd week (row 25) gets a start/end time pair
#    (5:30 AM - 9:45 AM)
# ------------------------------------------------------------------
$ws.Cells.Item(25, 3).Value = 0.2291666666666667
$ws.Cells.Item(25, 4).Value = 0.40625
$ws.Range("C25:D25").NumberFormat = "h:mm:ss"

# ------------------------------------------------------------------
# 4) The date column display now needs a timestamp component too.
#    Re-apply the (new) number format to every date cell across all
#    four week blocks so they all continue to share one style.
#    (Applied per contiguous block - this COM bridge only honours the
#    first area of a multi-area "A,B,C" range assignment.)
# ------------------------------------------------------------------
$ws.Range("B2:B8").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B11:B17").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B20:B26").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B29:B35").NumberFormat = "yyyy-mm-dd h:mm:ss"
